$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (header) ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Data table updates (rows 14-29) ---
# Row 14
$ws.Range("C14").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

# Row 15
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = "'***.*"
$ws.Range("E15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = -12.5
$ws.Range("N15").Value = 16.666666666666

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -35.294117647058
$ws.Range("I16").Value = 75
$ws.Range("J16").Value = 68
$ws.Range("K16").Value = 10.294117647058
$ws.Range("L16").Value = 36.363636363636
$ws.Range("M16").Value = -17.582417582417
$ws.Range("N16").Value = -70.238095238095

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 26.923076923076
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 110
$ws.Range("K17").Value = 13.636363636363
$ws.Range("L17").Value = 42.045454545454
$ws.Range("M17").Value = 54.320987654321
$ws.Range("N17").Value = 47.058823529411

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 140
$ws.Range("I18").Value = 82
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 82.222222222222
$ws.Range("L18").Value = 41.379310344827
$ws.Range("M18").Value = -19.607843137254
$ws.Range("N18").Value = -84.586466165413

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -27.777777777777
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -18.367346938775
$ws.Range("I19").Value = 193
$ws.Range("J19").Value = 228
$ws.Range("K19").Value = -15.350877192982
$ws.Range("L19").Value = 62.184873949579
$ws.Range("M19").Value = 77.064220183486
$ws.Range("N19").Value = 6.043956043956

# Row 20
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 150
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 90
$ws.Range("I20").Value = 146
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 21.666666666666
$ws.Range("L20").Value = 97.297297297297
$ws.Range("M20").Value = 94.666666666666
$ws.Range("N20").Value = -76.260162601626

# Row 21
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = 8.823529411764
$ws.Range("F21").Value = 150
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = 21.951219512195
$ws.Range("I21").Value = 630
$ws.Range("J21").Value = 579
$ws.Range("K21").Value = 8.808290155440
$ws.Range("L21").Value = 56.327543424317
$ws.Range("M21").Value = 34.615384615384
$ws.Range("N21").Value = -62.410501193317

# Row 22
$ws.Range("D22").Value = "'0"
$ws.Range("C22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "'***.*"
$ws.Range("N22").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0

# Row 23
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -40
$ws.Range("F23").Value = 14
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 45
$ws.Range("J23").Value = 36
$ws.Range("K23").Value = 25
$ws.Range("L23").Value = 80
$ws.Range("M23").Value = 125

# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -3.571428571428
$ws.Range("F24").Value = 137
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 39.795918367346
$ws.Range("I24").Value = 505
$ws.Range("J24").Value = 386
$ws.Range("K24").Value = 30.829015544041
$ws.Range("L24").Value = 53.030303030303
$ws.Range("M24").Value = 79.715302491103

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 9
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 28.571428571428
$ws.Range("I25").Value = 182
$ws.Range("J25").Value = 155
$ws.Range("K25").Value = 17.419354838709
$ws.Range("L25").Value = 61.061946902654
$ws.Range("M25").Value = -11.219512195122

# Row 26
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 12
$ws.Range("K26").Value = 9.090909090909
$ws.Range("L26").Value = 9.090909090909

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -71.428571428571
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 38.461538461538

# Row 28
$ws.Range("C28").Value = "'0"
$ws.Range("A28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = "'0"
$ws.Range("A28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "'***.*"
$ws.Range("A28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -88.888888888888
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = -68.75

# Row 29
$ws.Range("C29").Value = "'0"
$ws.Range("A29").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = "'0"
$ws.Range("A29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "'***.*"
$ws.Range("A29").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -80
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = -61.538461538461

$excel.CutCopyMode = $false
